# "Fitur Add Pelanggan dan Import Excel Done"
# Update the sample/template row: the package/plan name placeholder was
# changed from "inflyajalah" to "infly".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G6").Value = "infly"

# Reflect the cursor/selection position left behind by the editing session
# (scrolled right a bit, ending with M16 selected).
$ws.Range("M16").Select()
